$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.178.59"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.071.20"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.14"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.60%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "2.376.69"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.75"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.753"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "2.070.40"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").Value = "38.109.61"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.07"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +6.01%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.01%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.66%  "
$ws.Range("E38").Value = "  +4.49%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").Value = "1.483.14"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0951"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.81"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +14.51%  "
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "2.260.94"
$ws.Range("E51").Value = "  +2.36%  "
